# "select credential type to run"
#
# Inserts a new "CREDENTIAL_TYPE" row (with a "window" value and a
# window/asset list-validation) just above the existing "1_state" section
# header in the base sheet, shifting every later row down by one, and
# resizes the Tabla1 table + sheet dimension to match the extra row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row above row 22 ("1_state" section header) -------
$ws.Rows(22).Insert()

# --- 2. Populate the new row ------------------------------------------
$ws.Range("A22").Value = "CREDENTIAL_TYPE"
$ws.Range("C22").Value = "window"

# Vertically center the new row's cells (matches the rest of the sheet).
$ws.Range("A22:E22").VerticalAlignment = -4108

# --- 3. Data validation: dropdown list "window,asset" on C22 ----------
$ws.Range("C22").Validation.Add(3, 1, 1, '"window,asset"')

# --- 4. Grow the table (Tabla1) + autofilter by one row ----------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E74"))

# --- 5. Restore view state: selection on the new "1_state" area -------
$ws.Range("A27").Select()
